$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample-csv-products")

# Update D2:D148 to value 5
$ws.Range("D2:D148").Value = 5

# Select D2:D148 with active cell D2, matching the saved selection in the diff
$ws.Activate()
$ws.Range("D2:D148").Select()
